$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("D28").Value = 44859
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 4500
$ws.Range("M28").Value = 4250
$ws.Range("N28").Value = '$/paquete'
$ws.Range("O28").Value = 'Región de Arica y Parinacota'
$ws.Range("P28").Value = 4250
$ws.Range("Q28").Value = 1

# Row 29
$ws.Range("D29").Value = 44460
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 4000
$ws.Range("L29").Value = 4500
$ws.Range("M29").Value = 4250
$ws.Range("N29").Value = '$/paquete'
$ws.Range("O29").Value = 'Región de Arica y Parinacota'
$ws.Range("P29").Value = 4250
$ws.Range("Q29").Value = 1

# Row 30
$ws.Range("D30").Value = 44715
$ws.Range("J30").Value = 1600
$ws.Range("K30").Value = 3500
$ws.Range("L30").Value = 4000
$ws.Range("M30").Value = 3750
$ws.Range("N30").Value = '$/paquete'
$ws.Range("O30").Value = 'Región de Arica y Parinacota'
$ws.Range("P30").Value = 3750
$ws.Range("Q30").Value = 1

# Row 31
$ws.Range("D31").Value = 44698
$ws.Range("J31").Value = 1600
$ws.Range("K31").Value = 3800
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = 3900
$ws.Range("N31").Value = '$/paquete'
$ws.Range("O31").Value = 'Región de Arica y Parinacota'
$ws.Range("P31").Value = 3900
$ws.Range("Q31").Value = 1

# Row 32
$ws.Range("D32").Value = 44804
$ws.Range("J32").Value = 1200
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 4500
$ws.Range("M32").Value = 4250
$ws.Range("N32").Value = '$/paquete'
$ws.Range("O32").Value = 'Región de Arica y Parinacota'
$ws.Range("P32").Value = 4250
$ws.Range("Q32").Value = 1

# Row 33
$ws.Range("D33").Value = 44761
$ws.Range("J33").Value = 1400
$ws.Range("K33").Value = 4000
$ws.Range("L33").Value = 4500
$ws.Range("M33").Value = 4250
$ws.Range("N33").Value = '$/paquete'
$ws.Range("O33").Value = 'Región de Arica y Parinacota'
$ws.Range("P33").Value = 4250
$ws.Range("Q33").Value = 1

# Row 34
$ws.Range("D34").Value = 44476
$ws.Range("J34").Value = 600
$ws.Range("K34").Value = 3500
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = 3750
$ws.Range("N34").Value = '$/paquete'
$ws.Range("O34").Value = 'Región de Arica y Parinacota'
$ws.Range("P34").Value = 3750
$ws.Range("Q34").Value = 1

# Row 35
$ws.Range("D35").Value = 44445
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 4500
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = 4750
$ws.Range("N35").Value = '$/paquete'
$ws.Range("O35").Value = 'Región de Arica y Parinacota'
$ws.Range("P35").Value = 4750
$ws.Range("Q35").Value = 1

# Row 36
$ws.Range("D36").Value = 44515
$ws.Range("J36").Value = 800
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = 3500
$ws.Range("N36").Value = '$/paquete'
$ws.Range("O36").Value = 'Región de Arica y Parinacota'
$ws.Range("P36").Value = 3500
$ws.Range("Q36").Value = 1

# Row 37
$ws.Range("D37").Value = 44790
$ws.Range("J37").Value = 1600
$ws.Range("K37").Value = 3300
$ws.Range("L37").Value = 3500
$ws.Range("M37").Value = 3400
$ws.Range("N37").Value = '$/paquete'
$ws.Range("O37").Value = 'Región de Arica y Parinacota'
$ws.Range("P37").Value = 3400
$ws.Range("Q37").Value = 1

# Row 38
$ws.Range("D38").Value = 44418
$ws.Range("J38").Value = 800
$ws.Range("K38").Value = 4500
$ws.Range("L38").Value = 5000
$ws.Range("M38").Value = 4750
$ws.Range("N38").Value = '$/paquete'
$ws.Range("O38").Value = 'Región de Arica y Parinacota'
$ws.Range("P38").Value = 4750
$ws.Range("Q38").Value = 1

# Row 39
$ws.Range("D39").Value = 44841
$ws.Range("J39").Value = 1200
$ws.Range("K39").Value = 4000
$ws.Range("L39").Value = 4500
$ws.Range("M39").Value = 4250
$ws.Range("N39").Value = '$/paquete'
$ws.Range("O39").Value = 'Región de Arica y Parinacota'
$ws.Range("P39").Value = 4250
$ws.Range("Q39").Value = 1

# Row 40
$ws.Range("D40").Value = 44778
$ws.Range("J40").Value = 1120
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = 3750
$ws.Range("N40").Value = '$/paquete'
$ws.Range("O40").Value = 'Región de Arica y Parinacota'
$ws.Range("P40").Value = 3750
$ws.Range("Q40").Value = 1

# Row 41
$ws.Range("D41").Value = 44811
$ws.Range("J41").Value = 1200
$ws.Range("K41").Value = 4000
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = 4250
$ws.Range("N41").Value = '$/paquete'
$ws.Range("O41").Value = 'Región de Arica y Parinacota'
$ws.Range("P41").Value = 4250
$ws.Range("Q41").Value = 1

# Row 42
$ws.Range("D42").Value = 44446
$ws.Range("J42").Value = 800
$ws.Range("K42").Value = 4500
$ws.Range("L42").Value = 5000
$ws.Range("M42").Value = 4750
$ws.Range("N42").Value = '$/paquete'
$ws.Range("O42").Value = 'Región de Arica y Parinacota'
$ws.Range("P42").Value = 4750
$ws.Range("Q42").Value = 1

# Row 43
$ws.Range("D43").Value = 44631
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = 7500
$ws.Range("N43").Value = '$/docena de matas'
$ws.Range("O43").Value = 'Provincia del Elquí'
$ws.Range("P43").Value = 1250
$ws.Range("Q43").Value = 6

# Row 44
$ws.Range("D44").Value = 44473
$ws.Range("J44").Value = 600
$ws.Range("K44").Value = 3500
$ws.Range("L44").Value = 4000
$ws.Range("M44").Value = 3750
$ws.Range("N44").Value = '$/paquete'
$ws.Range("O44").Value = 'Región de Arica y Parinacota'
$ws.Range("P44").Value = 3750
$ws.Range("Q44").Value = 1

# Row 45
$ws.Range("D45").Value = 44512
$ws.Range("J45").Value = 800
$ws.Range("K45").Value = 3000
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = 3250
$ws.Range("N45").Value = '$/paquete'
$ws.Range("O45").Value = 'Región de Arica y Parinacota'
$ws.Range("P45").Value = 3250
$ws.Range("Q45").Value = 1

# Row 46
$ws.Range("D46").Value = 44532
$ws.Range("J46").Value = 740
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = 3500
$ws.Range("N46").Value = '$/paquete'
$ws.Range("O46").Value = 'Región de Arica y Parinacota'
$ws.Range("P46").Value = 3500
$ws.Range("Q46").Value = 1

# Row 47
$ws.Range("D47").Value = 44673
$ws.Range("J47").Value = 800
$ws.Range("K47").Value = 5000
$ws.Range("L47").Value = 5500
$ws.Range("M47").Value = 5250
$ws.Range("N47").Value = '$/docena de matas'
$ws.Range("O47").Value = 'Provincia del Elquí'
$ws.Range("P47").Value = 875
$ws.Range("Q47").Value = 6

# Row 48
$ws.Range("D48").Value = 44315
$ws.Range("J48").Value = 700
$ws.Range("K48").Value = 2500
$ws.Range("L48").Value = 3000
$ws.Range("M48").Value = 2750
$ws.Range("N48").Value = '$/paquete'
$ws.Range("O48").Value = 'Región de Arica y Parinacota'
$ws.Range("P48").Value = 2750
$ws.Range("Q48").Value = 1

# Row 49
$ws.Range("D49").Value = 44509
$ws.Range("J49").Value = 800
$ws.Range("K49").Value = 3500
$ws.Range("L49").Value = 4000
$ws.Range("M49").Value = 3750
$ws.Range("N49").Value = '$/paquete'
$ws.Range("O49").Value = 'Región de Arica y Parinacota'
$ws.Range("P49").Value = 3750
$ws.Range("Q49").Value = 1

# Row 50
$ws.Range("D50").Value = 44342
$ws.Range("J50").Value = 560
$ws.Range("K50").Value = 3000
$ws.Range("L50").Value = 3500
$ws.Range("M50").Value = 3250
$ws.Range("N50").Value = '$/paquete'
$ws.Range("O50").Value = 'Región de Arica y Parinacota'
$ws.Range("P50").Value = 3250
$ws.Range("Q50").Value = 1

# Row 51
$ws.Range("D51").Value = 44530
$ws.Range("J51").Value = 800
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 4000
$ws.Range("M51").Value = 3500
$ws.Range("N51").Value = '$/paquete'
$ws.Range("O51").Value = 'Región de Arica y Parinacota'
$ws.Range("P51").Value = 3500
$ws.Range("Q51").Value = 1

# Row 52
$ws.Range("D52").Value = 44819
$ws.Range("J52").Value = 1200
$ws.Range("K52").Value = 3800
$ws.Range("L52").Value = 4000
$ws.Range("M52").Value = 3900
$ws.Range("N52").Value = '$/paquete'
$ws.Range("O52").Value = 'Región de Arica y Parinacota'
$ws.Range("P52").Value = 3900
$ws.Range("Q52").Value = 1

# Row 53
$ws.Range("D53").Value = 44397
$ws.Range("J53").Value = 800
$ws.Range("K53").Value = 4000
$ws.Range("L53").Value = 4500
$ws.Range("M53").Value = 4250
$ws.Range("N53").Value = '$/paquete'
$ws.Range("O53").Value = 'Región de Arica y Parinacota'
$ws.Range("P53").Value = 4250
$ws.Range("Q53").Value = 1

# Row 54
$ws.Range("D54").Value = 44474
$ws.Range("J54").Value = 760
$ws.Range("K54").Value = 3500
$ws.Range("L54").Value = 4000
$ws.Range("M54").Value = 3750
$ws.Range("N54").Value = '$/paquete'
$ws.Range("O54").Value = 'Región de Arica y Parinacota'
$ws.Range("P54").Value = 3750
$ws.Range("Q54").Value = 1

# Row 55
$ws.Range("D55").Value = 44546
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 5000
$ws.Range("L55").Value = 5500
$ws.Range("M55").Value = 5250
$ws.Range("N55").Value = '$/docena de matas'
$ws.Range("O55").Value = 'Provincia del Elquí'
$ws.Range("P55").Value = 875
$ws.Range("Q55").Value = 6

# Row 56
$ws.Range("D56").Value = 44162
$ws.Range("J56").Value = 2000
$ws.Range("K56").Value = 2800
$ws.Range("L56").Value = 3000
$ws.Range("M56").Value = 2900
$ws.Range("N56").Value = '$/paquete'
$ws.Range("O56").Value = 'Región de Arica y Parinacota'
$ws.Range("P56").Value = 2900
$ws.Range("Q56").Value = 1

# Row 57
$ws.Range("D57").Value = 44365
$ws.Range("J57").Value = 800
$ws.Range("K57").Value = 3500
$ws.Range("L57").Value = 4000
$ws.Range("M57").Value = 3750
$ws.Range("N57").Value = '$/paquete'
$ws.Range("O57").Value = 'Región de Arica y Parinacota'
$ws.Range("P57").Value = 3750
$ws.Range("Q57").Value = 1

# Row 58
$ws.Range("D58").Value = 44680
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 5500
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = 5750
$ws.Range("N58").Value = '$/docena de matas'
$ws.Range("O58").Value = 'Provincia del Elquí'
$ws.Range("P58").Value = 958
$ws.Range("Q58").Value = 6

# Row 59
$ws.Range("D59").Value = 44525
$ws.Range("J59").Value = 720
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 4000
$ws.Range("M59").Value = 3500
$ws.Range("N59").Value = '$/paquete'
$ws.Range("O59").Value = 'Región de Arica y Parinacota'
$ws.Range("P59").Value = 3500
$ws.Range("Q59").Value = 1

# Row 60
$ws.Range("D60").Value = 44536
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 3500
$ws.Range("L60").Value = 4000
$ws.Range("M60").Value = 3750
$ws.Range("N60").Value = '$/paquete'
$ws.Range("O60").Value = 'Región de Arica y Parinacota'
$ws.Range("P60").Value = 3750
$ws.Range("Q60").Value = 1

# Row 61
$ws.Range("D61").Value = 44750
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 3500
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = 3750
$ws.Range("N61").Value = '$/paquete'
$ws.Range("O61").Value = 'Región de Arica y Parinacota'
$ws.Range("P61").Value = 3750
$ws.Range("Q61").Value = 1

# Row 62
$ws.Range("D62").Value = 44719
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = 3750
$ws.Range("N62").Value = '$/paquete'
$ws.Range("O62").Value = 'Región de Arica y Parinacota'
$ws.Range("P62").Value = 3750
$ws.Range("Q62").Value = 1

# Row 63
$ws.Range("D63").Value = 44453
$ws.Range("J63").Value = 800
$ws.Range("K63").Value = 4500
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = 4750
$ws.Range("N63").Value = '$/paquete'
$ws.Range("O63").Value = 'Región de Arica y Parinacota'
$ws.Range("P63").Value = 4750
$ws.Range("Q63").Value = 1

# Row 64
$ws.Range("D64").Value = 44356
$ws.Range("J64").Value = 600
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = 3250
$ws.Range("N64").Value = '$/paquete'
$ws.Range("O64").Value = 'Región de Arica y Parinacota'
$ws.Range("P64").Value = 3250
$ws.Range("Q64").Value = 1

# Row 65
$ws.Range("D65").Value = 44540
$ws.Range("J65").Value = 500
$ws.Range("K65").Value = 3000
$ws.Range("L65").Value = 4000
$ws.Range("M65").Value = 3500
$ws.Range("N65").Value = '$/paquete'
$ws.Range("O65").Value = 'Región de Arica y Parinacota'
$ws.Range("P65").Value = 3500
$ws.Range("Q65").Value = 1

# Row 66
$ws.Range("D66").Value = 44694
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 5000
$ws.Range("L66").Value = 5500
$ws.Range("M66").Value = 5250
$ws.Range("N66").Value = '$/docena de matas'
$ws.Range("O66").Value = 'Provincia del Elquí'
$ws.Range("P66").Value = 875
$ws.Range("Q66").Value = 6

# Row 67
$ws.Range("D67").Value = 44553
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 3500
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = 3750
$ws.Range("N67").Value = '$/docena de matas'
$ws.Range("O67").Value = 'Provincia del Elquí'
$ws.Range("P67").Value = 625
$ws.Range("Q67").Value = 6

# Row 68
$ws.Range("D68").Value = 44831
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 4500
$ws.Range("M68").Value = 4250
$ws.Range("N68").Value = '$/paquete'
$ws.Range("O68").Value = 'Región de Arica y Parinacota'
$ws.Range("P68").Value = 4250
$ws.Range("Q68").Value = 1

# Row 69
$ws.Range("D69").Value = 44701
$ws.Range("J69").Value = 1120
$ws.Range("K69").Value = 4000
$ws.Range("L69").Value = 4500
$ws.Range("M69").Value = 4250
$ws.Range("N69").Value = '$/paquete'
$ws.Range("O69").Value = 'Región de Arica y Parinacota'
$ws.Range("P69").Value = 4250
$ws.Range("Q69").Value = 1

# Row 70
$ws.Range("D70").Value = 44165
$ws.Range("J70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 3500
$ws.Range("M70").Value = 3250
$ws.Range("N70").Value = '$/paquete'
$ws.Range("O70").Value = 'Región de Arica y Parinacota'
$ws.Range("P70").Value = 3250
$ws.Range("Q70").Value = 1

# Row 71
$ws.Range("D71").Value = 44519
$ws.Range("J71").Value = 800
$ws.Range("K71").Value = 3500
$ws.Range("L71").Value = 4000
$ws.Range("M71").Value = 3750
$ws.Range("N71").Value = '$/paquete'
$ws.Range("O71").Value = 'Región de Arica y Parinacota'
$ws.Range("P71").Value = 3750
$ws.Range("Q71").Value = 1

# Row 72
$ws.Range("D72").Value = 44411
$ws.Range("J72").Value = 880
$ws.Range("K72").Value = 4000
$ws.Range("L72").Value = 4500
$ws.Range("M72").Value = 4250
$ws.Range("N72").Value = '$/paquete'
$ws.Range("O72").Value = 'Región de Arica y Parinacota'
$ws.Range("P72").Value = 4250
$ws.Range("Q72").Value = 1

# Row 73
$ws.Range("D73").Value = 44817
$ws.Range("J73").Value = 1200
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = 4250
$ws.Range("N73").Value = '$/paquete'
$ws.Range("O73").Value = 'Región de Arica y Parinacota'
$ws.Range("P73").Value = 4250
$ws.Range("Q73").Value = 1

# Row 74
$ws.Range("D74").Value = 44533
$ws.Range("J74").Value = 900
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = 3500
$ws.Range("N74").Value = '$/paquete'
$ws.Range("O74").Value = 'Región de Arica y Parinacota'
$ws.Range("P74").Value = 3500
$ws.Range("Q74").Value = 1

# Row 75
$ws.Range("D75").Value = 44516
$ws.Range("J75").Value = 740
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 4000
$ws.Range("M75").Value = 3500
$ws.Range("N75").Value = '$/paquete'
$ws.Range("O75").Value = 'Región de Arica y Parinacota'
$ws.Range("P75").Value = 3500
$ws.Range("Q75").Value = 1

# Row 76
$ws.Range("D76").Value = 44448
$ws.Range("J76").Value = 640
$ws.Range("K76").Value = 4500
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = 4750
$ws.Range("N76").Value = '$/paquete'
$ws.Range("O76").Value = 'Región de Arica y Parinacota'
$ws.Range("P76").Value = 4750
$ws.Range("Q76").Value = 1

# Row 77
$ws.Range("D77").Value = 44469
$ws.Range("J77").Value = 700
$ws.Range("K77").Value = 4000
$ws.Range("L77").Value = 4500
$ws.Range("M77").Value = 4250
$ws.Range("N77").Value = '$/paquete'
$ws.Range("O77").Value = 'Región de Arica y Parinacota'
$ws.Range("P77").Value = 4250
$ws.Range("Q77").Value = 1

# Row 78
$ws.Range("D78").Value = 44505
$ws.Range("J78").Value = 800
$ws.Range("K78").Value = 3500
$ws.Range("L78").Value = 4000
$ws.Range("M78").Value = 3750
$ws.Range("N78").Value = '$/paquete'
$ws.Range("O78").Value = 'Región de Arica y Parinacota'
$ws.Range("P78").Value = 3750
$ws.Range("Q78").Value = 1

# Row 79
$ws.Range("D79").Value = 44449
$ws.Range("J79").Value = 700
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = 4250
$ws.Range("N79").Value = '$/paquete'
$ws.Range("O79").Value = 'Región de Arica y Parinacota'
$ws.Range("P79").Value = 4250
$ws.Range("Q79").Value = 1

# Row 80
$ws.Range("D80").Value = 44176
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = 3250
$ws.Range("N80").Value = '$/paquete'
$ws.Range("O80").Value = 'Región de Arica y Parinacota'
$ws.Range("P80").Value = 3250
$ws.Range("Q80").Value = 1

# Row 81
$ws.Range("D81").Value = 44848
$ws.Range("J81").Value = 1100
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = 4250
$ws.Range("N81").Value = '$/paquete'
$ws.Range("O81").Value = 'Región de Arica y Parinacota'
$ws.Range("P81").Value = 4250
$ws.Range("Q81").Value = 1

# Row 82
$ws.Range("D82").Value = 44172
$ws.Range("J82").Value = 760
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 3500
$ws.Range("M82").Value = 3250
$ws.Range("N82").Value = '$/paquete'
$ws.Range("O82").Value = 'Región de Arica y Parinacota'
$ws.Range("P82").Value = 3250
$ws.Range("Q82").Value = 1

# Row 83
$ws.Range("D83").Value = 44452
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 5000
$ws.Range("M83").Value = 4750
$ws.Range("N83").Value = '$/paquete'
$ws.Range("O83").Value = 'Región de Arica y Parinacota'
$ws.Range("P83").Value = 4750
$ws.Range("Q83").Value = 1

# Row 84
$ws.Range("D84").Value = 44435
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 4500
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = 4750
$ws.Range("N84").Value = '$/paquete'
$ws.Range("O84").Value = 'Región de Arica y Parinacota'
$ws.Range("P84").Value = 4750
$ws.Range("Q84").Value = 1

# Row 85
$ws.Range("D85").Value = 44349
$ws.Range("J85").Value = 560
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 3500
$ws.Range("M85").Value = 3250
$ws.Range("N85").Value = '$/paquete'
$ws.Range("O85").Value = 'Región de Arica y Parinacota'
$ws.Range("P85").Value = 3250
$ws.Range("Q85").Value = 1

# Row 86
$ws.Range("D86").Value = 44522
$ws.Range("J86").Value = 800
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = 3250
$ws.Range("N86").Value = '$/paquete'
$ws.Range("O86").Value = 'Región de Arica y Parinacota'
$ws.Range("P86").Value = 3250
$ws.Range("Q86").Value = 1

# Row 87
$ws.Range("D87").Value = 44586
$ws.Range("J87").Value = 760
$ws.Range("K87").Value = 3500
$ws.Range("L87").Value = 4000
$ws.Range("M87").Value = 3750
$ws.Range("N87").Value = '$/docena de matas'
$ws.Range("O87").Value = 'Provincia del Elquí'
$ws.Range("P87").Value = 625
$ws.Range("Q87").Value = 6

# Row 88
$ws.Range("D88").Value = 44537
$ws.Range("J88").Value = 760
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = 3750
$ws.Range("N88").Value = '$/paquete'
$ws.Range("O88").Value = 'Región de Arica y Parinacota'
$ws.Range("P88").Value = 3750
$ws.Range("Q88").Value = 1

# Row 89
$ws.Range("D89").Value = 44818
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 4000
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 4250
$ws.Range("N89").Value = '$/paquete'
$ws.Range("O89").Value = 'Región de Arica y Parinacota'
$ws.Range("P89").Value = 4250
$ws.Range("Q89").Value = 1

# Row 90
$ws.Range("D90").Value = 44806
$ws.Range("J90").Value = 1000
$ws.Range("K90").Value = 4000
$ws.Range("L90").Value = 4500
$ws.Range("M90").Value = 4250
$ws.Range("N90").Value = '$/paquete'
$ws.Range("O90").Value = 'Región de Arica y Parinacota'
$ws.Range("P90").Value = 4250
$ws.Range("Q90").Value = 1

# Row 91
$ws.Range("D91").Value = 44748
$ws.Range("J91").Value = 1000
$ws.Range("K91").Value = 3500
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = 3750
$ws.Range("N91").Value = '$/paquete'
$ws.Range("O91").Value = 'Región de Arica y Parinacota'
$ws.Range("P91").Value = 3750
$ws.Range("Q91").Value = 1

# Row 92
$ws.Range("D92").Value = 44581
$ws.Range("J92").Value = 760
$ws.Range("K92").Value = 3500
$ws.Range("L92").Value = 4000
$ws.Range("M92").Value = 3750
$ws.Range("N92").Value = '$/docena de matas'
$ws.Range("O92").Value = 'Provincia del Elquí'
$ws.Range("P92").Value = 625
$ws.Range("Q92").Value = 6

# Row 93
$ws.Range("D93").Value = 44771
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 3500
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = 3750
$ws.Range("N93").Value = '$/paquete'
$ws.Range("O93").Value = 'Región de Arica y Parinacota'
$ws.Range("P93").Value = 3750
$ws.Range("Q93").Value = 1

# Row 94
$ws.Range("D94").Value = 44783
$ws.Range("J94").Value = 1600
$ws.Range("K94").Value = 3300
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = 3400
$ws.Range("N94").Value = '$/paquete'
$ws.Range("O94").Value = 'Región de Arica y Parinacota'
$ws.Range("P94").Value = 3400
$ws.Range("Q94").Value = 1

# Row 95
$ws.Range("D95").Value = 44810
$ws.Range("J95").Value = 1000
$ws.Range("K95").Value = 4000
$ws.Range("L95").Value = 4500
$ws.Range("M95").Value = 4250
$ws.Range("N95").Value = '$/paquete'
$ws.Range("O95").Value = 'Región de Arica y Parinacota'
$ws.Range("P95").Value = 4250
$ws.Range("Q95").Value = 1

# Row 96
$ws.Range("D96").Value = 44526
$ws.Range("J96").Value = 800
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 4000
$ws.Range("M96").Value = 3500
$ws.Range("N96").Value = '$/paquete'
$ws.Range("O96").Value = 'Región de Arica y Parinacota'
$ws.Range("P96").Value = 3500
$ws.Range("Q96").Value = 1

# Row 97
$ws.Range("D97").Value = 44582
$ws.Range("J97").Value = 640
$ws.Range("K97").Value = 3500
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = 3750
$ws.Range("N97").Value = '$/docena de matas'
$ws.Range("O97").Value = 'Provincia del Elquí'
$ws.Range("P97").Value = 625
$ws.Range("Q97").Value = 6

# Row 98
$ws.Range("D98").Value = 44377
$ws.Range("J98").Value = 600
$ws.Range("K98").Value = 4000
$ws.Range("L98").Value = 4500
$ws.Range("M98").Value = 4250
$ws.Range("N98").Value = '$/paquete'
$ws.Range("O98").Value = 'Región de Arica y Parinacota'
$ws.Range("P98").Value = 4250
$ws.Range("Q98").Value = 1

# Row 99
$ws.Range("D99").Value = 44504
$ws.Range("J99").Value = 760
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = 3250
$ws.Range("N99").Value = '$/paquete'
$ws.Range("O99").Value = 'Región de Arica y Parinacota'
$ws.Range("P99").Value = 3250
$ws.Range("Q99").Value = 1

# Row 100
$ws.Range("D100").Value = 44670
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 5000
$ws.Range("L100").Value = 5500
$ws.Range("M100").Value = 5250
$ws.Range("N100").Value = '$/docena de matas'
$ws.Range("O100").Value = 'Provincia del Elquí'
$ws.Range("P100").Value = 875
$ws.Range("Q100").Value = 6

# Row 101
$ws.Range("D101").Value = 44760
$ws.Range("J101").Value = 1200
$ws.Range("K101").Value = 3500
$ws.Range("L101").Value = 4000
$ws.Range("M101").Value = 3750
$ws.Range("N101").Value = '$/paquete'
$ws.Range("O101").Value = 'Región de Arica y Parinacota'
$ws.Range("P101").Value = 3750
$ws.Range("Q101").Value = 1

# Row 102
$ws.Range("D102").Value = 44769
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 3300
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 3400
$ws.Range("N102").Value = '$/paquete'
$ws.Range("O102").Value = 'Región de Arica y Parinacota'
$ws.Range("P102").Value = 3400
$ws.Range("Q102").Value = 1

# Row 103
$ws.Range("D103").Value = 44855
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 4000
$ws.Range("L103").Value = 4500
$ws.Range("M103").Value = 4250
$ws.Range("N103").Value = '$/paquete'
$ws.Range("O103").Value = 'Región de Arica y Parinacota'
$ws.Range("P103").Value = 4250
$ws.Range("Q103").Value = 1

# Row 104
$ws.Range("D104").Value = 44334
$ws.Range("J104").Value = 760
$ws.Range("K104").Value = 3000
$ws.Range("L104").Value = 3500
$ws.Range("M104").Value = 3250
$ws.Range("N104").Value = '$/paquete'
$ws.Range("O104").Value = 'Región de Arica y Parinacota'
$ws.Range("P104").Value = 3250
$ws.Range("Q104").Value = 1

# Row 105
$ws.Range("D105").Value = 44427
$ws.Range("J105").Value = 600
$ws.Range("K105").Value = 4500
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = 4750
$ws.Range("N105").Value = '$/paquete'
$ws.Range("O105").Value = 'Región de Arica y Parinacota'
$ws.Range("P105").Value = 4750
$ws.Range("Q105").Value = 1

# Row 106
$ws.Range("D106").Value = 44539
$ws.Range("J106").Value = 600
$ws.Range("K106").Value = 3000
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = 3500
$ws.Range("N106").Value = '$/paquete'
$ws.Range("O106").Value = 'Región de Arica y Parinacota'
$ws.Range("P106").Value = 3500
$ws.Range("Q106").Value = 1

# Row 107
$ws.Range("D107").Value = 44518
$ws.Range("J107").Value = 760
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = 3500
$ws.Range("N107").Value = '$/paquete'
$ws.Range("O107").Value = 'Región de Arica y Parinacota'
$ws.Range("P107").Value = 3500
$ws.Range("Q107").Value = 1

# Row 108
$ws.Range("D108").Value = 44488
$ws.Range("J108").Value = 800
$ws.Range("K108").Value = 3500
$ws.Range("L108").Value = 4000
$ws.Range("M108").Value = 3750
$ws.Range("N108").Value = '$/paquete'
$ws.Range("O108").Value = 'Región de Arica y Parinacota'
$ws.Range("P108").Value = 3750
$ws.Range("Q108").Value = 1

# Row 109
$ws.Range("D109").Value = 44348
$ws.Range("J109").Value = 700
$ws.Range("K109").Value = 3000
$ws.Range("L109").Value = 3500
$ws.Range("M109").Value = 3250
$ws.Range("N109").Value = '$/paquete'
$ws.Range("O109").Value = 'Región de Arica y Parinacota'
$ws.Range("P109").Value = 3250
$ws.Range("Q109").Value = 1

# Row 110
$ws.Range("D110").Value = 44341
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = 3250
$ws.Range("N110").Value = '$/paquete'
$ws.Range("O110").Value = 'Región de Arica y Parinacota'
$ws.Range("P110").Value = 3250
$ws.Range("Q110").Value = 1

# Row 111
$ws.Range("D111").Value = 44169
$ws.Range("J111").Value = 2400
$ws.Range("K111").Value = 3000
$ws.Range("L111").Value = 3500
$ws.Range("M111").Value = 3250
$ws.Range("N111").Value = '$/paquete'
$ws.Range("O111").Value = 'Región de Arica y Parinacota'
$ws.Range("P111").Value = 3250
$ws.Range("Q111").Value = 1

# Row 112
$ws.Range("D112").Value = 44729
$ws.Range("J112").Value = 1140
$ws.Range("K112").Value = 3500
$ws.Range("L112").Value = 4000
$ws.Range("M112").Value = 3750
$ws.Range("N112").Value = '$/paquete'
$ws.Range("O112").Value = 'Región de Arica y Parinacota'
$ws.Range("P112").Value = 3750
$ws.Range("Q112").Value = 1

# Row 113
$ws.Range("D113").Value = 44463
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 4250
$ws.Range("N113").Value = '$/paquete'
$ws.Range("O113").Value = 'Región de Arica y Parinacota'
$ws.Range("P113").Value = 4250
$ws.Range("Q113").Value = 1

# Row 114
$ws.Range("D114").Value = 44685
$ws.Range("J114").Value = 2000
$ws.Range("K114").Value = 5000
$ws.Range("L114").Value = 5500
$ws.Range("M114").Value = 5250
$ws.Range("N114").Value = '$/docena de matas'
$ws.Range("O114").Value = 'Provincia del Elquí'
$ws.Range("P114").Value = 875
$ws.Range("Q114").Value = 6

# Row 115
$ws.Range("D115").Value = 44462
$ws.Range("J115").Value = 660
$ws.Range("K115").Value = 4000
$ws.Range("L115").Value = 4500
$ws.Range("M115").Value = 4250
$ws.Range("N115").Value = '$/paquete'
$ws.Range("O115").Value = 'Región de Arica y Parinacota'
$ws.Range("P115").Value = 4250
$ws.Range("Q115").Value = 1

# Row 116
$ws.Range("D116").Value = 44498
$ws.Range("J116").Value = 900
$ws.Range("K116").Value = 3800
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 3900
$ws.Range("N116").Value = '$/paquete'
$ws.Range("O116").Value = 'Región de Arica y Parinacota'
$ws.Range("P116").Value = 3900
$ws.Range("Q116").Value = 1

# Row 117
$ws.Range("D117").Value = 44455
$ws.Range("J117").Value = 600
$ws.Range("K117").Value = 4500
$ws.Range("L117").Value = 5000
$ws.Range("M117").Value = 4750
$ws.Range("N117").Value = '$/paquete'
$ws.Range("O117").Value = 'Región de Arica y Parinacota'
$ws.Range("P117").Value = 4750
$ws.Range("Q117").Value = 1

# Row 118
$ws.Range("D118").Value = 44484
$ws.Range("J118").Value = 840
$ws.Range("K118").Value = 3500
$ws.Range("L118").Value = 4000
$ws.Range("M118").Value = 3750
$ws.Range("N118").Value = '$/paquete'
$ws.Range("O118").Value = 'Región de Arica y Parinacota'
$ws.Range("P118").Value = 3750
$ws.Range("Q118").Value = 1

# Row 119
$ws.Range("D119").Value = 44434
$ws.Range("J119").Value = 600
$ws.Range("K119").Value = 4500
$ws.Range("L119").Value = 5000
$ws.Range("M119").Value = 4750
$ws.Range("N119").Value = '$/paquete'
$ws.Range("O119").Value = 'Región de Arica y Parinacota'
$ws.Range("P119").Value = 4750
$ws.Range("Q119").Value = 1

# Row 120
$ws.Range("D120").Value = 44441
$ws.Range("J120").Value = 600
$ws.Range("K120").Value = 4500
$ws.Range("L120").Value = 5000
$ws.Range("M120").Value = 4750
$ws.Range("N120").Value = '$/paquete'
$ws.Range("O120").Value = 'Región de Arica y Parinacota'
$ws.Range("P120").Value = 4750
$ws.Range("Q120").Value = 1

# Row 121
$ws.Range("D121").Value = 44799
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 4000
$ws.Range("L121").Value = 4500
$ws.Range("M121").Value = 4250
$ws.Range("N121").Value = '$/paquete'
$ws.Range("O121").Value = 'Región de Arica y Parinacota'
$ws.Range("P121").Value = 4250
$ws.Range("Q121").Value = 1

# Row 122
$ws.Range("D122").Value = 44736
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3500
$ws.Range("L122").Value = 4000
$ws.Range("M122").Value = 3750
$ws.Range("N122").Value = '$/paquete'
$ws.Range("O122").Value = 'Región de Arica y Parinacota'
$ws.Range("P122").Value = 3750
$ws.Range("Q122").Value = 1

# Row 123
$ws.Range("D123").Value = 44379
$ws.Range("J123").Value = 800
$ws.Range("K123").Value = 4000
$ws.Range("L123").Value = 4500
$ws.Range("M123").Value = 4250
$ws.Range("N123").Value = '$/paquete'
$ws.Range("O123").Value = 'Región de Arica y Parinacota'
$ws.Range("P123").Value = 4250
$ws.Range("Q123").Value = 1

# Row 124
$ws.Range("D124").Value = 44813
$ws.Range("J124").Value = 1300
$ws.Range("K124").Value = 4000
$ws.Range("L124").Value = 4500
$ws.Range("M124").Value = 4250
$ws.Range("N124").Value = '$/paquete'
$ws.Range("O124").Value = 'Región de Arica y Parinacota'
$ws.Range("P124").Value = 4250
$ws.Range("Q124").Value = 1

# Row 125
$ws.Range("D125").Value = 44832
$ws.Range("J125").Value = 1400
$ws.Range("K125").Value = 4000
$ws.Range("L125").Value = 4500
$ws.Range("M125").Value = 4250
$ws.Range("N125").Value = '$/paquete'
$ws.Range("O125").Value = 'Región de Arica y Parinacota'
$ws.Range("P125").Value = 4250
$ws.Range("Q125").Value = 1

# Row 126
$ws.Range("D126").Value = 44722
$ws.Range("J126").Value = 1100
$ws.Range("K126").Value = 3500
$ws.Range("L126").Value = 4000
$ws.Range("M126").Value = 3750
$ws.Range("N126").Value = '$/paquete'
$ws.Range("O126").Value = 'Región de Arica y Parinacota'
$ws.Range("P126").Value = 3750
$ws.Range("Q126").Value = 1

# Row 127
$ws.Range("D127").Value = 44490
$ws.Range("J127").Value = 660
$ws.Range("K127").Value = 3500
$ws.Range("L127").Value = 4000
$ws.Range("M127").Value = 3750
$ws.Range("N127").Value = '$/paquete'
$ws.Range("O127").Value = 'Región de Arica y Parinacota'
$ws.Range("P127").Value = 3750
$ws.Range("Q127").Value = 1

# Row 128
$ws.Range("D128").Value = 44845
$ws.Range("J128").Value = 1200
$ws.Range("K128").Value = 4000
$ws.Range("L128").Value = 4500
$ws.Range("M128").Value = 4250
$ws.Range("N128").Value = '$/paquete'
$ws.Range("O128").Value = 'Región de Arica y Parinacota'
$ws.Range("P128").Value = 4250
$ws.Range("Q128").Value = 1

# Row 129
$ws.Range("D129").Value = 44497
$ws.Range("J129").Value = 740
$ws.Range("K129").Value = 3800
$ws.Range("L129").Value = 4000
$ws.Range("M129").Value = 3900
$ws.Range("N129").Value = '$/paquete'
$ws.Range("O129").Value = 'Región de Arica y Parinacota'
$ws.Range("P129").Value = 3900
$ws.Range("Q129").Value = 1

# Row 130
$ws.Range("A130").Value = 8
$ws.Range("B130").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C130").Value = 'Coquimbo'
$ws.Range("D130").Value = 44335
$ws.Range("J130").Value = 600
$ws.Range("K130").Value = 3000
$ws.Range("L130").Value = 3500
$ws.Range("M130").Value = 3250
$ws.Range("N130").Value = '$/paquete'
$ws.Range("O130").Value = 'Región de Arica y Parinacota'
$ws.Range("P130").Value = 3250
$ws.Range("Q130").Value = 1
$ws.Range("E130").Value = 4
$ws.Range("F130").Value = 100112052
$ws.Range("G130").Value = 'Albahaca'
$ws.Range("H130").Value = 'Sin especificar'
$ws.Range("I130").Value = 'Primera'
$ws.Range("R130").Value = 'Hortaliza'
$ws.Range("D130").NumberFormat = $ws.Range("D129").NumberFormat
